$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.461
$ws.Range("D3").Value = -7.623
$ws.Range("E19").Value = 16.515
$ws.Range("B21").Value = 9.458
$ws.Range("B23").Value = 7.891999999999999
$ws.Range("D24").Value = -7.279000000000001
$ws.Range("E24").Value = 16.652
$ws.Range("B25").Value = 6.396000000000001
$ws.Range("C27").Value = -13.055
$ws.Range("E30").Value = 16.531
$ws.Range("C31").Value = -12.85
$ws.Range("E31").Value = 16.318
$ws.Range("E33").Value = 17.435
$ws.Range("C39").Value = -12.847
$ws.Range("C48").Value = -11.094
$ws.Range("C51").Value = -11.329
$ws.Range("C52").Value = -11.389
$ws.Range("B53").Value = 5.895999999999999
$ws.Range("C55").Value = -13.656
$ws.Range("E55").Value = 16.441
$ws.Range("C56").Value = -12.996
$ws.Range("B57").Value = 5.230000000000001
$ws.Range("C57").Value = -13.342
$ws.Range("D57").Value = -8.279
$ws.Range("B59").Value = 4.714
$ws.Range("D61").Value = -7.739
$ws.Range("E65").Value = 17.267
$ws.Range("B69").Value = 5.667000000000001
$ws.Range("D70").Value = -7.390000000000001
$ws.Range("E70").Value = 17.659
$ws.Range("C73").Value = -12.818
$ws.Range("E75").Value = 16.738
$ws.Range("B79").Value = 5.76
$ws.Range("B83").Value = 5.542
$ws.Range("E83").Value = 16.885
$ws.Range("D86").Value = -8.186
$ws.Range("C89").Value = -11.273
$ws.Range("C90").Value = -12.714
$ws.Range("B93").Value = 5.659000000000001
$ws.Range("E96").Value = 16.454
$ws.Range("E97").Value = 17.16
$ws.Range("D98").Value = -8.397
$ws.Range("D100").Value = -8.361999999999998
$ws.Range("D102").Value = -7.805000000000001

$wb.Save()
